# Auto-generated edit script applying numeric updates from the commit diff
# to the "Masamune_Profits" workbook (8 item-category sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (51 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 51.82
$ws.Range("I15").Value = 51.82
$ws.Range("K15").Value = 155.46
$ws.Range("M15").Value = 13.53999999999999
$ws.Range("H19").Value = 643.6429000000001
$ws.Range("I19").Value = 723.44446
$ws.Range("J19").Value = 500
$ws.Range("K19").Value = 723.44446
$ws.Range("L19").Value = 500
$ws.Range("M19").Value = -548.44446
$ws.Range("N19").Value = -850
$ws.Range("H32").Value = 1283.4706
$ws.Range("I32").Value = 1111.375
$ws.Range("J32").Value = 1336.4231
$ws.Range("K32").Value = 1111.375
$ws.Range("L32").Value = 1336.4231
$ws.Range("M32").Value = -785.375
$ws.Range("N32").Value = -1988.4231
$ws.Range("H53").Value = 81.333336
$ws.Range("I53").Value = 53.153847
$ws.Range("J53").Value = 114.63636
$ws.Range("K53").Value = 53.153847
$ws.Range("L53").Value = 114.63636
$ws.Range("M53").Value = 583.846153
$ws.Range("N53").Value = -1388.63636
$ws.Range("H62").Value = 3705.3572
$ws.Range("I62").Value = 3759.6155
$ws.Range("K62").Value = 3759.6155
$ws.Range("M62").Value = -3135.6155
$ws.Range("H65").Value = 3705.3572
$ws.Range("I65").Value = 3759.6155
$ws.Range("K65").Value = 18798.0775
$ws.Range("M65").Value = -15678.0775
$ws.Range("H129").Value = 402779.9
$ws.Range("I129").Value = 707135.6
$ws.Range("J129").Value = 4776.231
$ws.Range("K129").Value = 2121406.8
$ws.Range("L129").Value = 14328.693
$ws.Range("M129").Value = -2116406.8
$ws.Range("N129").Value = -24328.693
$ws.Range("H134").Value = 73050
$ws.Range("J134").Value = 73050
$ws.Range("L134").Value = 73050
$ws.Range("N134").Value = -83190
$ws.Range("H138").Value = 3924.7473
$ws.Range("I138").Value = 3840.9167
$ws.Range("J138").Value = 3936.8674
$ws.Range("K138").Value = 11522.7501
$ws.Range("L138").Value = 11810.6022
$ws.Range("M138").Value = -6382.750100000001
$ws.Range("N138").Value = -22090.6022

# --- Sheet: ARM (35 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29705.56
$ws.Range("I32").Value = 30036.46
$ws.Range("J32").Value = 29132
$ws.Range("K32").Value = 30036.46
$ws.Range("L32").Value = 29132
$ws.Range("M32").Value = -29749.46
$ws.Range("N32").Value = -29706
$ws.Range("H74").Value = 2711.889
$ws.Range("I74").Value = 1992.0952
$ws.Range("J74").Value = 5231.1665
$ws.Range("K74").Value = 1992.0952
$ws.Range("L74").Value = 5231.1665
$ws.Range("M74").Value = -1118.0952
$ws.Range("N74").Value = -6979.1665
$ws.Range("H77").Value = 2711.889
$ws.Range("I77").Value = 1992.0952
$ws.Range("J77").Value = 5231.1665
$ws.Range("K77").Value = 9960.476000000001
$ws.Range("L77").Value = 26155.8325
$ws.Range("M77").Value = -5592.476000000001
$ws.Range("N77").Value = -34891.8325
$ws.Range("H102").Value = 12172.685
$ws.Range("I102").Value = 1646.6666
$ws.Range("J102").Value = 30217.285
$ws.Range("K102").Value = 1646.6666
$ws.Range("L102").Value = 30217.285
$ws.Range("M102").Value = -24.66660000000002
$ws.Range("N102").Value = -33461.285
$ws.Range("H110").Value = 1946.2727
$ws.Range("I110").Value = 1837.7333
$ws.Range("J110").Value = 2178.8572
$ws.Range("K110").Value = 1837.7333
$ws.Range("L110").Value = 2178.8572
$ws.Range("M110").Value = 207.2666999999999
$ws.Range("N110").Value = -6268.8572

# --- Sheet: BSM (25 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 413
$ws.Range("I64").Value = 458.2
$ws.Range("J64").Value = 380.7143
$ws.Range("K64").Value = 458.2
$ws.Range("L64").Value = 380.7143
$ws.Range("M64").Value = -233.2
$ws.Range("N64").Value = -830.7143
$ws.Range("H67").Value = 413
$ws.Range("I67").Value = 458.2
$ws.Range("J67").Value = 380.7143
$ws.Range("K67").Value = 458.2
$ws.Range("L67").Value = 380.7143
$ws.Range("M67").Value = 321.8
$ws.Range("N67").Value = -1940.7143
$ws.Range("H107").Value = 2108.5386
$ws.Range("I107").Value = 1899.7142
$ws.Range("J107").Value = 2352.1667
$ws.Range("K107").Value = 1899.7142
$ws.Range("L107").Value = 2352.1667
$ws.Range("M107").Value = 20.28580000000011
$ws.Range("N107").Value = -6192.1667
$ws.Range("H139").Value = 73333
$ws.Range("J139").Value = 73333
$ws.Range("L139").Value = 73333
$ws.Range("N139").Value = -83613

# --- Sheet: CRP (23 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H81").Value = 49924
$ws.Range("J81").Value = 49924
$ws.Range("L81").Value = 49924
$ws.Range("N81").Value = -51920
$ws.Range("H84").Value = 49924
$ws.Range("J84").Value = 49924
$ws.Range("L84").Value = 149772
$ws.Range("N84").Value = -159756
$ws.Range("H105").Value = 2345.1875
$ws.Range("I105").Value = 2345.1875
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2345.1875
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -598.1875
$ws.Range("N105").ClearContents()
$ws.Range("H131").Value = 35326
$ws.Range("J131").Value = 35326
$ws.Range("L131").Value = 35326
$ws.Range("N131").Value = -45406
$ws.Range("H137").Value = 66314.875
$ws.Range("J137").Value = 66314.875
$ws.Range("L137").Value = 66314.875
$ws.Range("N137").Value = -76514.875

# --- Sheet: CUL (37 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 249.05882
$ws.Range("I23").Value = 207
$ws.Range("J23").Value = 258.07144
$ws.Range("K23").Value = 621
$ws.Range("L23").Value = 774.21432
$ws.Range("M23").Value = -386
$ws.Range("N23").Value = -1244.21432
$ws.Range("H49").Value = 3000
$ws.Range("J49").Value = 3000
$ws.Range("L49").Value = 9000
$ws.Range("N49").Value = -9312
$ws.Range("H88").Value = 7013.1763
$ws.Range("J88").Value = 7013.1763
$ws.Range("L88").Value = 21039.5289
$ws.Range("N88").Value = -21895.5289
$ws.Range("H91").Value = 7013.1763
$ws.Range("J91").Value = 7013.1763
$ws.Range("L91").Value = 21039.5289
$ws.Range("N91").Value = -24003.5289
$ws.Range("H112").Value = 3308.6296
$ws.Range("I112").Value = 799.6667
$ws.Range("J112").Value = 3405.1282
$ws.Range("K112").Value = 2399.0001
$ws.Range("L112").Value = 10215.3846
$ws.Range("M112").Value = -1291.0001
$ws.Range("N112").Value = -12431.3846
$ws.Range("H130").Value = 188951.44
$ws.Range("J130").Value = 1654.0769
$ws.Range("L130").Value = 4962.2307
$ws.Range("N130").Value = -15002.2307
$ws.Range("H131").Value = 903.8
$ws.Range("I131").Value = 498.22223
$ws.Range("J131").Value = 943.9121
$ws.Range("K131").Value = 1494.66669
$ws.Range("L131").Value = 2831.7363
$ws.Range("M131").Value = 3545.33331
$ws.Range("N131").Value = -12911.7363

# --- Sheet: GSM (15 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2592.3333
$ws.Range("I113").Value = 2222
$ws.Range("J113").Value = 3333
$ws.Range("K113").Value = 2222
$ws.Range("L113").Value = 3333
$ws.Range("M113").Value = -52
$ws.Range("N113").Value = -7673
$ws.Range("H120").Value = 39317
$ws.Range("J120").Value = 39317
$ws.Range("L120").Value = 39317
$ws.Range("N120").Value = -48993
$ws.Range("H127").Value = 65884
$ws.Range("J127").Value = 65884
$ws.Range("L127").Value = 65884
$ws.Range("N127").Value = -75804

# --- Sheet: LTW (32 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 70800
$ws.Range("J76").Value = 70800
$ws.Range("L76").Value = 70800
$ws.Range("N76").Value = -71476
$ws.Range("H79").Value = 70800
$ws.Range("J79").Value = 70800
$ws.Range("L79").Value = 70800
$ws.Range("N79").Value = -73140
$ws.Range("H92").Value = 23964
$ws.Range("J92").Value = 23964
$ws.Range("L92").Value = 23964
$ws.Range("N92").Value = -28956
$ws.Range("H109").Value = 35281
$ws.Range("J109").Value = 35281
$ws.Range("L109").Value = 35281
$ws.Range("N109").Value = -38055
$ws.Range("H117").Value = 45392
$ws.Range("J117").Value = 45392
$ws.Range("L117").Value = 45392
$ws.Range("N117").Value = -54570
$ws.Range("H123").Value = 42429
$ws.Range("J123").Value = 42429
$ws.Range("L123").Value = 42429
$ws.Range("N123").Value = -52229
$ws.Range("H131").Value = 33826
$ws.Range("J131").Value = 33826
$ws.Range("L131").Value = 33826
$ws.Range("N131").Value = -43906
$ws.Range("H133").Value = 35292
$ws.Range("J133").Value = 35292
$ws.Range("L133").Value = 35292
$ws.Range("N133").Value = -40352

# --- Sheet: WVR (23 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 31592.857
$ws.Range("J93").Value = 31592.857
$ws.Range("L93").Value = 31592.857
$ws.Range("N93").Value = -36584.857
$ws.Range("H107").Value = 20000238
$ws.Range("I107").Value = 263.66666
$ws.Range("J107").Value = 50000200
$ws.Range("K107").Value = 790.9999799999999
$ws.Range("L107").Value = 150000600
$ws.Range("M107").Value = 1129.00002
$ws.Range("N107").Value = -150004440
$ws.Range("H118").Value = 43388
$ws.Range("J118").Value = 43388
$ws.Range("L118").Value = 43388
$ws.Range("N118").Value = -46702
$ws.Range("H127").Value = 36657.332
$ws.Range("J127").Value = 36657.332
$ws.Range("L127").Value = 36657.332
$ws.Range("N127").Value = -46577.332
$ws.Range("H139").Value = 56875
$ws.Range("J139").Value = 56875
$ws.Range("L139").Value = 56875
$ws.Range("N139").Value = -67155

